# Update NATMI LR-pair values with new TPM-derived results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5668126666666667
$ws.Range("H2").Value = 1.700438
$ws.Range("I2").Value = 0.8744630508533847
$ws.Range("J2").Value = 0.8744630508533846
$ws.Range("M2").Value = 10.56874366666667
$ws.Range("N2").Value = 31.706231
$ws.Range("O2").Value = 0.04029387683847273
$ws.Range("P2").Value = 0.04029387683847273
$ws.Range("Q2").Value = 5.990497781019778
$ws.Range("R2").Value = 53.914480029178
$ws.Range("S2").Value = 0.0352355064708814
$ws.Range("T2").Value = 0.0352355064708814
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5668126666666667
$ws.Range("H3").Value = 1.700438
$ws.Range("I3").Value = 0.8744630508533847
$ws.Range("J3").Value = 0.8744630508533846
$ws.Range("M3").Value = 13.69721566666667
$ws.Range("N3").Value = 41.091647
$ws.Range("O3").Value = 0.05222133666117545
$ws.Range("P3").Value = 0.05222133666117545
$ws.Range("Q3").Value = 7.763755337931779
$ws.Range("R3").Value = 69.87379804138601
$ws.Range("S3").Value = 0.04566562937637319
$ws.Range("T3").Value = 0.04566562937637319
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5668126666666667
$ws.Range("H4").Value = 1.700438
$ws.Range("I4").Value = 0.8744630508533847
$ws.Range("J4").Value = 0.8744630508533846
$ws.Range("M4").Value = 71.48326
$ws.Range("N4").Value = 214.44978
$ws.Range("O4").Value = 0.2725335920045991
$ws.Range("P4").Value = 0.2725335920045991
$ws.Range("Q4").Value = 40.51761722262668
$ws.Range("R4").Value = 364.65855500364
$ws.Range("S4").Value = 0.2383205563243734
$ws.Range("T4").Value = 0.2383205563243734
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5668126666666667
$ws.Range("H5").Value = 1.700438
$ws.Range("I5").Value = 0.8744630508533847
$ws.Range("J5").Value = 0.8744630508533846
$ws.Range("M5").Value = 7.094812
$ws.Range("N5").Value = 21.284436
$ws.Range("O5").Value = 0.02704933433306391
$ws.Range("P5").Value = 0.02704933433306391
$ws.Range("Q5").Value = 4.021429309218667
$ws.Range("R5").Value = 36.192863782968
$ws.Range("S5").Value = 0.02365364342444427
$ws.Range("T5").Value = 0.02365364342444427
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5668126666666667
$ws.Range("H6").Value = 1.700438
$ws.Range("I6").Value = 0.8744630508533847
$ws.Range("J6").Value = 0.8744630508533846
$ws.Range("M6").Value = 88.21463766666666
$ws.Range("N6").Value = 264.643913
$ws.Range("O6").Value = 0.3363228267804361
$ws.Range("P6").Value = 0.3363228267804361
$ws.Range("Q6").Value = 50.00117401487712
$ws.Range("R6").Value = 450.010566133894
$ws.Range("S6").Value = 0.2941018851780546
$ws.Range("T6").Value = 0.2941018851780545
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5668126666666667
$ws.Range("H7").Value = 1.700438
$ws.Range("I7").Value = 0.8744630508533847
$ws.Range("J7").Value = 0.8744630508533846
$ws.Range("M7").Value = 71.23288733333334
$ws.Range("N7").Value = 213.698662
$ws.Range("O7").Value = 0.2715790333822526
$ws.Range("P7").Value = 0.2715790333822526
$ws.Range("Q7").Value = 40.3757028237729
$ws.Range("R7").Value = 363.3813254139561
$ws.Range("S7").Value = 0.2374858300792579
$ws.Range("T7").Value = 0.2374858300792578
$ws.Range("I8").Value = 0.07095005479414014
$ws.Range("J8").Value = 0.07095005479414014
$ws.Range("M8").Value = 10.56874366666667
$ws.Range("N8").Value = 31.706231
$ws.Range("O8").Value = 0.04029387683847273
$ws.Range("P8").Value = 0.04029387683847273
$ws.Range("Q8").Value = 0.4860424295717778
$ws.Range("R8").Value = 4.374381866146
$ws.Range("S8").Value = 0.002858852769557975
$ws.Range("T8").Value = 0.002858852769557975
$ws.Range("I9").Value = 0.07095005479414014
$ws.Range("J9").Value = 0.07095005479414014
$ws.Range("M9").Value = 13.69721566666667
$ws.Range("N9").Value = 41.091647
$ws.Range("O9").Value = 0.05222133666117545
$ws.Range("P9").Value = 0.05222133666117545
$ws.Range("Q9").Value = 0.6299166855557778
$ws.Range("R9").Value = 5.669250170002001
$ws.Range("S9").Value = 0.003705106697533637
$ws.Range("T9").Value = 0.003705106697533638
$ws.Range("I10").Value = 0.07095005479414014
$ws.Range("J10").Value = 0.07095005479414014
$ws.Range("M10").Value = 71.48326
$ws.Range("N10").Value = 214.44978
$ws.Range("O10").Value = 0.2725335920045991
$ws.Range("P10").Value = 0.2725335920045991
$ws.Range("Q10").Value = 3.287419816386667
$ws.Range("R10").Value = 29.58677834748
$ws.Range("S10").Value = 0.01933627328597014
$ws.Range("T10").Value = 0.01933627328597014
$ws.Range("I11").Value = 0.07095005479414014
$ws.Range("J11").Value = 0.07095005479414014
$ws.Range("M11").Value = 7.094812
$ws.Range("N11").Value = 21.284436
$ws.Range("O11").Value = 0.02704933433306391
$ws.Range("P11").Value = 0.02704933433306391
$ws.Range("Q11").Value = 0.3262809441306667
$ws.Range("R11").Value = 2.936528497176
$ws.Range("S11").Value = 0.0019191517530759
$ws.Range("T11").Value = 0.0019191517530759
$ws.Range("I12").Value = 0.07095005479414014
$ws.Range("J12").Value = 0.07095005479414014
$ws.Range("M12").Value = 88.21463766666666
$ws.Range("N12").Value = 264.643913
$ws.Range("O12").Value = 0.3363228267804361
$ws.Range("P12").Value = 0.3363228267804361
$ws.Range("Q12").Value = 4.056873566773112
$ws.Range("R12").Value = 36.511862100958
$ws.Range("S12").Value = 0.02386212298859204
$ws.Range("T12").Value = 0.02386212298859204
$ws.Range("I13").Value = 0.07095005479414014
$ws.Range("J13").Value = 0.07095005479414014
$ws.Range("M13").Value = 71.23288733333334
$ws.Range("N13").Value = 213.698662
$ws.Range("O13").Value = 0.2715790333822526
$ws.Range("P13").Value = 0.2715790333822526
$ws.Range("Q13").Value = 3.27590551127689
$ws.Range("R13").Value = 29.483149601492
$ws.Range("S13").Value = 0.01926854729941044
$ws.Range("T13").Value = 0.01926854729941044
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.03538233333333334
$ws.Range("H14").Value = 0.106147
$ws.Range("I14").Value = 0.0545868943524752
$ws.Range("J14").Value = 0.0545868943524752
$ws.Range("M14").Value = 10.56874366666667
$ws.Range("N14").Value = 31.706231
$ws.Range("O14").Value = 0.04029387683847273
$ws.Range("P14").Value = 0.04029387683847273
$ws.Range("Q14").Value = 0.3739468113285556
$ws.Range("R14").Value = 3.365521301957
$ws.Range("S14").Value = 0.002199517598033359
$ws.Range("T14").Value = 0.002199517598033359
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.03538233333333334
$ws.Range("H15").Value = 0.106147
$ws.Range("I15").Value = 0.0545868943524752
$ws.Range("J15").Value = 0.0545868943524752
$ws.Range("M15").Value = 13.69721566666667
$ws.Range("N15").Value = 41.091647
$ws.Range("O15").Value = 0.05222133666117545
$ws.Range("P15").Value = 0.05222133666117545
$ws.Range("Q15").Value = 0.4846394504565556
$ws.Range("R15").Value = 4.361755054109
$ws.Range("S15").Value = 0.002850600587268624
$ws.Range("T15").Value = 0.002850600587268624
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.03538233333333334
$ws.Range("H16").Value = 0.106147
$ws.Range("I16").Value = 0.0545868943524752
$ws.Range("J16").Value = 0.0545868943524752
$ws.Range("M16").Value = 71.48326
$ws.Range("N16").Value = 214.44978
$ws.Range("O16").Value = 0.2725335920045991
$ws.Range("P16").Value = 0.2725335920045991
$ws.Range("Q16").Value = 2.529244533073333
$ws.Range("R16").Value = 22.76320079766
$ws.Range("S16").Value = 0.01487676239425563
$ws.Range("T16").Value = 0.01487676239425563
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.03538233333333334
$ws.Range("H17").Value = 0.106147
$ws.Range("I17").Value = 0.0545868943524752
$ws.Range("J17").Value = 0.0545868943524752
$ws.Range("M17").Value = 7.094812
$ws.Range("N17").Value = 21.284436
$ws.Range("O17").Value = 0.02704933433306391
$ws.Range("P17").Value = 0.02704933433306391
$ws.Range("Q17").Value = 0.2510310031213334
$ws.Range("R17").Value = 2.259279028092
$ws.Range("S17").Value = 0.00147653915554374
$ws.Range("T17").Value = 0.00147653915554374
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.03538233333333334
$ws.Range("H18").Value = 0.106147
$ws.Range("I18").Value = 0.0545868943524752
$ws.Range("J18").Value = 0.0545868943524752
$ws.Range("M18").Value = 88.21463766666666
$ws.Range("N18").Value = 264.643913
$ws.Range("O18").Value = 0.3363228267804361
$ws.Range("P18").Value = 0.3363228267804361
$ws.Range("Q18").Value = 3.121239714801222
$ws.Range("R18").Value = 28.091157433211
$ws.Range("S18").Value = 0.01835881861378948
$ws.Range("T18").Value = 0.01835881861378948
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.03538233333333334
$ws.Range("H19").Value = 0.106147
$ws.Range("I19").Value = 0.0545868943524752
$ws.Range("J19").Value = 0.0545868943524752
$ws.Range("M19").Value = 71.23288733333334
$ws.Range("N19").Value = 213.698662
$ws.Range("O19").Value = 0.2715790333822526
$ws.Range("P19").Value = 0.2715790333822526
$ws.Range("Q19").Value = 2.520385763923778
$ws.Range("R19").Value = 22.683471875314
$ws.Range("S19").Value = 0.01482465600358436
$ws.Range("T19").Value = 0.01482465600358436
